$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")
$project = $wb.Worksheets.Item("project")

# Header row: new columns for groomed_file / alignment_file
$data.Range("B1").Value = "groomed_file"
$data.Range("B2").Value = "./groomed/ellipsoid_1_DT.nrrd"
$data.Range("C1").Value = "alignment_file"
$data.Range("C2").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -39.999849 -49.999811 -49.999811"

$data.Range("B3").Value = "./groomed/ellipsoid_2_DT.nrrd"
$data.Range("C3").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -42.499850 -50.000000 -50.000000"

$data.Range("B4").Value = "./groomed/ellipsoid_3_DT.nrrd"
$data.Range("C4").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -44.999848 -49.999831 -49.999831"

$data.Range("B5").Value = "./groomed/ellipsoid_4_DT.nrrd"
$data.Range("C5").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -47.499850 -50.000000 -50.000000"

$data.Range("B6").Value = "./groomed/ellipsoid_5_DT.nrrd"
$data.Range("C6").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -49.999851 -49.999851 -49.999851"

$data.Range("B7").Value = "./groomed/ellipsoid_6_DT.nrrd"
$data.Range("C7").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -52.499850 -50.000000 -50.000000"

$data.Range("B8").Value = "./groomed/ellipsoid_7_DT.nrrd"
$data.Range("C8").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -54.999699 -49.999862 -49.999862"

$data.Range("B9").Value = "./groomed/ellipsoid_8_DT.nrrd"
$data.Range("C9").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -57.499850 -50.000000 -50.000000"

$data.Range("B10").Value = "./groomed/ellipsoid_9_DT.nrrd"
$data.Range("C10").Value = " 1.000000 0.000000 0.000000 0.000000 1.000000 0.000000 0.000000 0.000000 1.000000 -59.999700 -49.999873 -49.999873"

# Bump project version from 1 to 2 (keep it stored as text, like the rest
# of the key/value sheets, instead of letting it coerce to a number)
$scratch = $project.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "2"
$scratch.Copy()
$project.Range("B2").PasteSpecial(-4163)
$scratch.Clear()

# Make the "data" sheet the active sheet/tab (was "studio")
$data.Activate()
